$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be stored as text so numeric-looking strings
# such as "0.3100" or "10.90" keep their exact original formatting,
# matching how the source data (inline text strings) is represented.
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '19.929.89'
$ws.Range('E2').Value = '  -5.93%  '
$ws.Range('D3').Value = '1.408.53'
$ws.Range('E3').Value = '  -7.13%  '
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  -0.57%  '
$ws.Range('E5').Value = '  -0.49%  '
$ws.Range('D6').Value = '276.09'
$ws.Range('E6').Value = '  -3.72%  '
$ws.Range('D7').Value = '0.3659'
$ws.Range('E7').Value = '  -5.85%  '
$ws.Range('D8').Value = '0.3100'
$ws.Range('E8').Value = '  -1.51%  '
$ws.Range('D9').Value = '39.64'
$ws.Range('E9').Value = '  -6.44%  '
$ws.Range('D10').Value = '1.033'
$ws.Range('E10').Value = '  -1.91%  '
$ws.Range('D11').Value = '0.06511'
$ws.Range('E11').Value = '  -7.50%  '
$ws.Range('D12').Value = '1.003'
$ws.Range('E12').Value = '  -0.63%  '
$ws.Range('D13').Value = '5.484'
$ws.Range('E13').Value = '  -3.05%  '
$ws.Range('D14').Value = '17.62'
$ws.Range('E14').Value = '  -1.91%  '
$ws.Range('D15').Value = '6.184'
$ws.Range('E15').Value = '  -3.05%  '
$ws.Range('D16').Value = '1.410.31'
$ws.Range('E16').Value = '  -7.46%  '
$ws.Range('D17').Value = '0.00001018'
$ws.Range('E17').Value = '  -5.63%  '
$ws.Range('D18').Value = '0.05671'
$ws.Range('E18').Value = '  -13.95%  '
$ws.Range('E19').Value = '  -0.56%  '
$ws.Range('D20').Value = '70.90'
$ws.Range('E20').Value = '  -13.85%  '
$ws.Range('D21').Value = '5.610'
$ws.Range('E21').Value = '  -6.91%  '
$ws.Range('D22').Value = '14.72'
$ws.Range('E22').Value = '  -3.34%  '
$ws.Range('D23').Value = '10.90'
$ws.Range('E23').Value = '  +1.61%  '
$ws.Range('D24').Value = '2.230'
$ws.Range('E24').Value = '  -5.29%  '
$ws.Range('D25').Value = '19.939.87'
$ws.Range('E25').Value = '  -5.93%  '
$ws.Range('D26').Value = '2.258'
$ws.Range('E26').Value = '  -4.71%  '
$ws.Range('D27').Value = '132.76'
$ws.Range('E27').Value = '  -9.97%  '
$ws.Range('D28').Value = '17.28'
$ws.Range('E28').Value = '  -5.09%  '
$ws.Range('D29').Value = '1.569.90'
$ws.Range('E29').Value = '  -7.25%  '
$ws.Range('D30').Value = '109.53'
$ws.Range('E30').Value = '  -5.38%  '
$ws.Range('D31').Value = '3.909'
$ws.Range('E31').Value = '  -18.69%  '
$ws.Range('D32').Value = '5.268'
$ws.Range('E32').Value = '  -11.92%  '
$ws.Range('D33').Value = '0.8151'
$ws.Range('E33').Value = '  -14.16%  '
$ws.Range('D34').Value = '0.07694'
$ws.Range('E34').Value = '  -3.68%  '
$ws.Range('B35').Value = 'WEMIXTOKEN'
$ws.Range('C35').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D35').Value = '1.478'
$ws.Range('E35').Value = '  -0.52%  '
$ws.Range('B36').Value = 'FraxShare'
$ws.Range('C36').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D36').Value = '8.311'
$ws.Range('E36').Value = '  -1.71%  '
$ws.Range('D37').Value = '4.918'
$ws.Range('E37').Value = '  -3.46%  '
$ws.Range('D38').Value = '0.05814'
$ws.Range('E38').Value = '  -0.60%  '
$ws.Range('E39').Value = '  -0.53%  '
$ws.Range('D40').Value = '0.02062'
$ws.Range('E40').Value = '  -4.05%  '
$ws.Range('D41').Value = '10.48'
$ws.Range('E41').Value = '  -7.37%  '
$ws.Range('D42').Value = '0.1892'
$ws.Range('E42').Value = '  -5.13%  '
$ws.Range('E43').Value = '  -6.10%  '
$ws.Range('D44').Value = '0.5308'
$ws.Range('E44').Value = '  -6.47%  '
$ws.Range('D45').Value = '12.33'
$ws.Range('E45').Value = '  -4.76%  '
$ws.Range('D46').Value = '3.540'
$ws.Range('E46').Value = '  -4.37%  '
$ws.Range('D47').Value = '0.5177'
$ws.Range('E47').Value = '  -5.64%  '
$ws.Range('D48').Value = '114.96'
$ws.Range('E48').Value = '  -0.02%  '
$ws.Range('D49').Value = '1.769'
$ws.Range('E49').Value = '  -5.04%  '
$ws.Range('D50').Value = '1.033'
$ws.Range('E50').Value = '  -9.52%  '
$ws.Range('D51').Value = '1.002'
$ws.Range('E51').Value = '  -0.58%  '
